$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4.981230423876468
$ws.Range("C2").Value = 1.533380357151145
$ws.Range("D2").Value = 0.0765853248849524
$ws.Range("E2").Value = 1.359268881961739
$ws.Range("F2").Value = 2.948547961839751
$ws.Range("G2").Value = 0.0007876981006996155
$ws.Range("H2").Value = 0.01175518264854031
$ws.Range("I2").Value = 0.006613395466705363
$ws.Range("P2").Value = 1.097255617216376
$ws.Range("B3").Value = 4.320021657265556
$ws.Range("C3").Value = 1.324722100298629
$ws.Range("D3").Value = 0.07095763349168038
$ws.Range("E3").Value = 1.172073800671214
$ws.Range("F3").Value = 2.613109361582872
$ws.Range("G3").Value = 0.0007946891518557333
$ws.Range("H3").Value = 0.007212128027621834
$ws.Range("I3").Value = 0.002948885422485059
$ws.Range("P3").Value = 1.094481175163864
$ws.Range("B4").Value = 3.915698898342782
$ws.Range("C4").Value = 1.198322027119559
$ws.Range("D4").Value = 0.06745113563161453
$ws.Range("E4").Value = 1.058569938732447
$ws.Range("F4").Value = 2.409144672161048
$ws.Range("G4").Value = 0.0007990927253246388
$ws.Range("H4").Value = 0.004916224201339192
$ws.Range("I4").Value = 0.00148838123403916
$ws.Range("P4").Value = 1.092854804039931
$ws.Range("B5").Value = 3.747607318179689
$ws.Range("C5").Value = 1.147871087881583
$ws.Range("D5").Value = 0.06584242137925855
$ws.Range("E5").Value = 1.012576996458506
$ws.Range("F5").Value = 2.323316459685529
$ws.Range("G5").Value = 0.0008009303667910995
$ws.Range("H5").Value = 0.00408997314903714
$ws.Range("I5").Value = 0.001131603897055466
$ws.Range("P5").Value = 1.09065326676032
$ws.Range("B6").Value = 3.715280378918862
$ws.Range("C6").Value = 1.140358903001413
$ws.Range("D6").Value = 0.06537209896362839
$ws.Range("E6").Value = 1.004904273159696
$ws.Range("F6").Value = 2.305276147110703
$ws.Range("G6").Value = 0.0008012544326233542
$ws.Range("H6").Value = 0.003955060321159998
$ws.Range("I6").Value = 0.00116115291958252
$ws.Range("P6").Value = 1.088405492363435
$ws.Range("B7").Value = 3.901266069596943
$ws.Range("C7").Value = 1.19995177266685
$ws.Range("D7").Value = 0.06687462319382576
$ws.Range("E7").Value = 1.05779948839654
$ws.Range("F7").Value = 2.397540478082703
$ws.Range("G7").Value = 0.0007991637887946879
$ws.Range("H7").Value = 0.004891676594446448
$ws.Range("I7").Value = 0.001675678934604186
$ws.Range("P7").Value = 1.087684391494015
$ws.Range("B8").Value = 4.736434825295305
$ws.Range("C8").Value = 1.464118339770437
$ws.Range("D8").Value = 0.07391235557678044
$ws.Range("E8").Value = 1.29416445724037
$ws.Range("F8").Value = 2.818429149018527
$ws.Range("G8").Value = 0.0007901460370847985
$ws.Range("H8").Value = 0.01005516924093632
$ws.Range("I8").Value = 0.005372860819191416
$ws.Range("P8").Value = 1.089487591704398
$ws.Range("B9").Value = 6.415911772227162
$ws.Range("C9").Value = 1.994740460635967
$ws.Range("D9").Value = 0.0883769866692532
$ws.Range("E9").Value = 1.772201791939992
$ws.Range("F9").Value = 3.683463379825696
$ws.Range("G9").Value = 0.0007732046513021398
$ws.Range("H9").Value = 0.02468761372549633
$ws.Range("I9").Value = 0.01955913258151032
$ws.Range("P9").Value = 1.103558314488197
$ws.Range("B10").Value = 7.602905785113592
$ws.Range("C10").Value = 2.376245661793007
$ws.Range("D10").Value = 0.09544131048739501
$ws.Range("E10").Value = 2.03044398709352
$ws.Range("F10").Value = 4.270257071634973
$ws.Range("G10").Value = 0.0007616360504728184
$ws.Range("H10").Value = 0.0378572013089622
$ws.Range("I10").Value = 0.03511095460282565
$ws.Range("P10").Value = 1.091024460534172
$ws.Range("B11").Value = 7.588332565320059
$ws.Range("C11").Value = 2.345199510257487
$ws.Range("D11").Value = 0.0728433864938367
$ws.Range("E11").Value = 1.330751556135809
$ws.Range("F11").Value = 3.982694589184518
$ws.Range("G11").Value = 0.0007601541995707229
$ws.Range("H11").Value = 0.05280914036437068
$ws.Range("I11").Value = 0.037753114894854
$ws.Range("P11").Value = 0.9030669903682167
$ws.Range("B12").Value = 7.345618155442139
$ws.Range("C12").Value = 2.23660020803942
$ws.Range("D12").Value = 0.05644071730884548
$ws.Range("E12").Value = 0.8211603905815821
$ws.Range("F12").Value = 3.643944556141406
$ws.Range("G12").Value = 0.0007609435587364955
$ws.Range("H12").Value = 0.08744102875231619
$ws.Range("I12").Value = 0.03644100085519231
$ws.Range("P12").Value = 0.7777018069335782
$ws.Range("B13").Value = 6.897314969269019
$ws.Range("C13").Value = 2.067507920854155
$ws.Range("D13").Value = 0.04305388161407464
$ws.Range("E13").Value = 0.4304368035855077
$ws.Range("F13").Value = 3.232595640091887
$ws.Range("G13").Value = 0.0007635703527378856
$ws.Range("H13").Value = 0.1385865610561012
$ws.Range("I13").Value = 0.03234772729667501
$ws.Range("P13").Value = 0.6847524645093088
$ws.Range("B14").Value = 6.49179278763927
$ws.Range("C14").Value = 1.92469935826449
$ws.Range("D14").Value = 0.03545943499494442
$ws.Range("E14").Value = 0.2304414329774218
$ws.Range("F14").Value = 2.914477541059597
$ws.Range("G14").Value = 0.0007662003112132409
$ws.Range("H14").Value = 0.1846997451284835
$ws.Range("I14").Value = 0.02857043298932282
$ws.Range("P14").Value = 0.636617747162866
$ws.Range("B15").Value = 6.345712463304778
$ws.Range("C15").Value = 1.87803954795487
$ws.Range("D15").Value = 0.03377085357345422
$ws.Range("E15").Value = 0.1903441791731133
$ws.Range("F15").Value = 2.817604012370936
$ws.Range("G15").Value = 0.0007672860252298415
$ws.Range("H15").Value = 0.1961007734092988
$ws.Range("I15").Value = 0.02717402711722627
$ws.Range("P15").Value = 0.6281378259490218
$ws.Range("B16").Value = 5.936119174746807
$ws.Range("C16").Value = 1.756462432841943
$ws.Range("D16").Value = 0.03384028511147363
$ws.Range("E16").Value = 0.1804204042892827
$ws.Range("F16").Value = 2.657050463300919
$ws.Range("G16").Value = 0.0007716125292273823
$ws.Range("H16").Value = 0.1794527841903033
$ws.Range("I16").Value = 0.02184367076104721
$ws.Range("P16").Value = 0.651435611291852
$ws.Range("B17").Value = 5.833874973497018
$ws.Range("C17").Value = 1.737849896047749
$ws.Range("D17").Value = 0.03750419642386404
$ws.Range("E17").Value = 0.2639330924608174
$ws.Range("F17").Value = 2.69897714991464
$ws.Range("G17").Value = 0.0007735824492763888
$ws.Range("H17").Value = 0.1402455617029688
$ws.Range("I17").Value = 0.01969550886613902
$ws.Range("P17").Value = 0.6925813884899128
$ws.Range("B18").Value = 5.994937761584595
$ws.Range("C18").Value = 1.804951456726997
$ws.Range("D18").Value = 0.04650595151861836
$ws.Range("E18").Value = 0.5012350193891635
$ws.Range("F18").Value = 2.933583967417889
$ws.Range("G18").Value = 0.000773548116443686
$ws.Range("H18").Value = 0.08838109454863741
$ws.Range("I18").Value = 0.01962906721980584
$ws.Range("P18").Value = 0.7668408410710015
$ws.Range("B19").Value = 6.340717363587203
$ws.Range("C19").Value = 1.94458517355838
$ws.Range("D19").Value = 0.06104873384579435
$ws.Range("E19").Value = 0.9453177939100215
$ws.Range("F19").Value = 3.302873377397788
$ws.Range("G19").Value = 0.0007717270255375275
$ws.Range("H19").Value = 0.04693725659048198
$ws.Range("I19").Value = 0.02186821644847559
$ws.Range("P19").Value = 0.8739235503081701
$ws.Range("B20").Value = 7.246826672728787
$ws.Range("C20").Value = 2.281277834259754
$ws.Range("D20").Value = 0.0917043484040363
$ws.Range("E20").Value = 1.956744648804346
$ws.Range("F20").Value = 4.078312994721784
$ws.Range("G20").Value = 0.0007648015064704452
$ws.Range("H20").Value = 0.03398296660081845
$ws.Range("I20").Value = 0.03094044142827457
$ws.Range("P20").Value = 1.077553176369335
$ws.Range("B21").Value = 8.271384052501162
$ws.Range("C21").Value = 2.616745136655311
$ws.Range("D21").Value = 0.1023510655173041
$ws.Range("E21").Value = 2.323912087395797
$ws.Range("F21").Value = 4.642254762457441
$ws.Range("G21").Value = 0.0007553058648042135
$ws.Range("H21").Value = 0.04725742275376454
$ws.Range("I21").Value = 0.0460517376191687
$ws.Range("P21").Value = 1.106713656422215
$ws.Range("B22").Value = 8.934333451542443
$ws.Range("C22").Value = 2.825058429238027
$ws.Range("D22").Value = 0.1084892713896934
$ws.Range("E22").Value = 2.514524627111911
$ws.Range("F22").Value = 4.999520173465328
$ws.Range("G22").Value = 0.0007493132598017731
$ws.Range("H22").Value = 0.05638144862008332
$ws.Range("I22").Value = 0.05706598819542208
$ws.Range("P22").Value = 1.121950448346524
$ws.Range("B23").Value = 8.595242952112244
$ws.Range("C23").Value = 2.7106406466051
$ws.Range("D23").Value = 0.1058939781486643
$ws.Range("E23").Value = 2.412708661291333
$ws.Range("F23").Value = 4.821158926806987
$ws.Range("G23").Value = 0.0007524584084325714
$ws.Range("H23").Value = 0.05146333044048212
$ws.Range("I23").Value = 0.05092019961416483
$ws.Range("P23").Value = 1.119745363144574
$ws.Range("B24").Value = 7.302328049739117
$ws.Range("C24").Value = 2.290898539287753
$ws.Range("D24").Value = 0.09503829538671482
$ws.Range("E24").Value = 2.034688665809696
$ws.Range("F24").Value = 4.138735064769037
$ws.Range("G24").Value = 0.0007645285058988202
$ws.Range("H24").Value = 0.03472078643161858
$ws.Range("I24").Value = 0.03097258645892609
$ws.Range("P24").Value = 1.103740619409947
$ws.Range("B25").Value = 5.935627322275309
$ws.Range("C25").Value = 1.853165442785212
$ws.Range("D25").Value = 0.08346568978059565
$ws.Range("E25").Value = 1.640566109694802
$ws.Range("F25").Value = 3.427107054559997
$ws.Range("G25").Value = 0.0007777875869371744
$ws.Range("H25").Value = 0.0201211202158289
$ws.Range("I25").Value = 0.01505530797214849
$ws.Range("P25").Value = 1.090272130385642
